$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(34).Insert()

$ws.Range("A34").Value = 11
$ws.Range("B34").Value = 'Vega Monumental Concepción'
$ws.Range("C34").Value = 'Bíobío'
$ws.Range("D34").Value = 44967
$ws.Range("E34").Value = 8
$ws.Range("F34").Value = 100112012
$ws.Range("G34").Value = 'Espinaca'
$ws.Range("H34").Value = 'Sin especificar'
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 40
$ws.Range("K34").Value = 7000
$ws.Range("L34").Value = 7500
$ws.Range("M34").Value = 7250
$ws.Range("N34").Value = '$/cuna 10 kilos'
$ws.Range("O34").Value = 'Región Metropolitana'
$ws.Range("P34").Value = 725
$ws.Range("Q34").Value = 10
$ws.Range("R34").Value = 'Hortaliza'
